# Update "想去人数" (number of people interested) values on the
# "展览" (Worksheets 1) and "全部类型" (Worksheets 4) sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3043
$ws1.Range("F4").Value = 46
$ws1.Range("F5").Value = 34
$ws1.Range("F7").Value = 217
$ws1.Range("F8").Value = 14553
$ws1.Range("F9").Value = 162
$ws1.Range("F11").Value = 5812
$ws1.Range("F12").Value = 590
$ws1.Range("F14").Value = 44
$ws1.Range("F15").Value = 64
$ws1.Range("F16").Value = 1238
$ws1.Range("F18").Value = 82
$ws1.Range("F19").Value = 184
$ws1.Range("F20").Value = 797
$ws1.Range("F22").Value = 56
$ws1.Range("F23").Value = 10601
$ws1.Range("F27").Value = 3736

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3043
$ws4.Range("F5").Value = 46
$ws4.Range("F6").Value = 34
$ws4.Range("F8").Value = 217
$ws4.Range("F9").Value = 14553
$ws4.Range("F10").Value = 162
$ws4.Range("F12").Value = 5812
$ws4.Range("F13").Value = 590
$ws4.Range("F15").Value = 44
$ws4.Range("F16").Value = 64
$ws4.Range("F17").Value = 1238
$ws4.Range("F19").Value = 82
$ws4.Range("F20").Value = 184
$ws4.Range("F21").Value = 797
$ws4.Range("F23").Value = 56
$ws4.Range("F25").Value = 10601
$ws4.Range("F29").Value = 3736
